$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.828.00'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '1.896.62'
$ws.Range("E3").Value = '  -0.89%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'0.7785"
$ws.Range("E5").Value = '  +5.22%  '
$ws.Range("D6").Value = "'240.14"
$ws.Range("E6").Value = '  -1.66%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'0.3065"
$ws.Range("E8").Value = '  -1.98%  '
$ws.Range("D9").Value = "'25.48"
$ws.Range("E9").Value = '  -5.13%  '
$ws.Range("D10").Value = "'0.06853"
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").Value = "'0.07986"
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = '1.913.29'
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = "'0.7369"
$ws.Range("E13").Value = '  -5.48%  '
$ws.Range("D14").Value = "'5.177"
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").Value = "'91.27"
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '29.841.09'
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'13.77"
$ws.Range("E17").Value = '  -4.58%  '
$ws.Range("D18").Value = "'5.883"
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = "'244.22"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").Value = "'0.000007700"
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D22").Value = '2.153.35'
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = "'6.917"
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'166.70"
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = "'9.272"
$ws.Range("E26").Value = '  -1.88%  '
$ws.Range("D27").Value = "'18.73"
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("D28").Value = "'0.1316"
$ws.Range("E28").Value = '  +2.22%  '
$ws.Range("D29").Value = "'2.024"
$ws.Range("E29").Value = '  -2.03%  '
$ws.Range("D30").Value = "'1.387"
$ws.Range("E30").Value = '  +2.38%  '
$ws.Range("D31").Value = "'1.509"
$ws.Range("E31").Value = '  -2.53%  '
$ws.Range("D32").Value = "'4.270"
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = "'4.062"
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").Value = "'0.05245"
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("D35").Value = "'1.241"
$ws.Range("E35").Value = '  -4.69%  '
$ws.Range("D36").Value = "'0.7276"
$ws.Range("E36").Value = '  -3.39%  '
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").Value = "'0.01904"
$ws.Range("E38").Value = '  -2.06%  '
$ws.Range("D39").Value = "'2.779"
$ws.Range("E39").Value = '  -0.83%  '
$ws.Range("D40").Value = "'6.178"
$ws.Range("E40").Value = '  -3.31%  '
$ws.Range("D41").Value = "'0.4415"
$ws.Range("E41").Value = '  -2.39%  '
$ws.Range("D42").Value = "'72.09"
$ws.Range("E42").Value = '  -4.22%  '
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D44").Value = "'0.8389"
$ws.Range("E44").Value = '  -0.08%  '
$ws.Range("D45").Value = "'1.878"
$ws.Range("E45").Value = '  -4.54%  '
$ws.Range("D46").Value = "'7.585"
$ws.Range("E46").Value = '  -3.35%  '
$ws.Range("D47").Value = "'100.19"
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("D48").Value = "'9.784"
$ws.Range("E48").Value = '  -1.62%  '
$ws.Range("D49").Value = '2.058.52'
$ws.Range("E49").Value = '  -0.27%  '
$ws.Range("D50").Value = "'36.07"
$ws.Range("E50").Value = '  -3.18%  '
$ws.Range("D51").Value = "'928.30"
$ws.Range("E51").Value = '  -1.02%  '
